# Apply the "Updated cryptos list" refresh: coin list shifted down one row
# (new "OKB" entry inserted at row 9) and prices/24h deltas updated to the
# latest scrape. Values are written as literal text (leading apostrophe)
# and the style is reset to "Normal" afterwards so price strings like
# "305.67" or "19.30" are not silently coerced into numbers and lose their
# original formatting, matching how the source sheet stores them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.527.85' }
    @{ Cell = 'E2'; Value = '  -0.05%  ' }
    @{ Cell = 'D3'; Value = '1.812.96' }
    @{ Cell = 'E3'; Value = '  +0.02%  ' }
    @{ Cell = 'E4'; Value = '  -0.47%  ' }
    @{ Cell = 'D6'; Value = '305.67' }
    @{ Cell = 'E6'; Value = '  -0.98%  ' }
    @{ Cell = 'E7'; Value = '  -0.45%  ' }
    @{ Cell = 'D8'; Value = '0.3593' }
    @{ Cell = 'E8'; Value = '  -2.03%  ' }
    @{ Cell = 'B9'; Value = 'OKB' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D9'; Value = '46.39' }
    @{ Cell = 'E9'; Value = '  +2.98%  ' }
    @{ Cell = 'B10'; Value = 'Dogecoin' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge' }
    @{ Cell = 'D10'; Value = '0.07118' }
    @{ Cell = 'E10'; Value = '  -0.23%  ' }
    @{ Cell = 'B11'; Value = 'Polygon' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D11'; Value = '0.8922' }
    @{ Cell = 'E11'; Value = '  +1.30%  ' }
    @{ Cell = 'B12'; Value = 'TRON' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' }
    @{ Cell = 'D12'; Value = '0.07704' }
    @{ Cell = 'E12'; Value = '  -0.80%  ' }
    @{ Cell = 'B13'; Value = 'Solana' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' }
    @{ Cell = 'D13'; Value = '19.30' }
    @{ Cell = 'E13'; Value = '  -0.41%  ' }
    @{ Cell = 'B14'; Value = 'WrappedEther' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D14'; Value = '1.803.42' }
    @{ Cell = 'E14'; Value = '  -0.66%  ' }
    @{ Cell = 'B15'; Value = 'Polkadot' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D15'; Value = '5.259' }
    @{ Cell = 'E15'; Value = '  -0.69%  ' }
    @{ Cell = 'B16'; Value = 'Chainlink' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = 'D16'; Value = '6.306' }
    @{ Cell = 'E16'; Value = '  -1.13%  ' }
    @{ Cell = 'B17'; Value = 'Litecoin' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D17'; Value = '86.02' }
    @{ Cell = 'E17'; Value = '  -0.77%  ' }
    @{ Cell = 'B18'; Value = 'BinanceUSD' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' }
    @{ Cell = 'D18'; Value = '1.005' }
    @{ Cell = 'E18'; Value = '  -0.43%  ' }
    @{ Cell = 'B19'; Value = 'ShibaInu' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D19'; Value = '0.000008546' }
    @{ Cell = 'E19'; Value = '  -0.55%  ' }
    @{ Cell = 'B20'; Value = 'Dai' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D20'; Value = '1.003' }
    @{ Cell = 'E20'; Value = '  -0.46%  ' }
    @{ Cell = 'B21'; Value = 'WrappedBTC' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = 'D21'; Value = '26.554.79' }
    @{ Cell = 'E21'; Value = '  -0.18%  ' }
    @{ Cell = 'B22'; Value = 'Avalanche' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'D22'; Value = '14.15' }
    @{ Cell = 'E22'; Value = '  -0.74%  ' }
    @{ Cell = 'B23'; Value = 'Uniswap' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D23'; Value = '4.963' }
    @{ Cell = 'E23'; Value = '  -0.97%  ' }
    @{ Cell = 'B24'; Value = 'Cosmos' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D24'; Value = '10.52' }
    @{ Cell = 'E24'; Value = '  +0.36%  ' }
    @{ Cell = 'B25'; Value = 'Toncoin' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D25'; Value = '1.925' }
    @{ Cell = 'E25'; Value = '  -2.95%  ' }
    @{ Cell = 'B26'; Value = 'Monero' }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D26'; Value = '151.76' }
    @{ Cell = 'E26'; Value = '  +0.21%  ' }
    @{ Cell = 'B27'; Value = 'EthereumClassic' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D27'; Value = '17.77' }
    @{ Cell = 'E27'; Value = '  -0.97%  ' }
    @{ Cell = 'B28'; Value = 'LidoDAOToken' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D28'; Value = '2.020' }
    @{ Cell = 'E28'; Value = '  -1.93%  ' }
    @{ Cell = 'B29'; Value = 'BitcoinCash' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D29'; Value = '112.13' }
    @{ Cell = 'E29'; Value = '  -0.86%  ' }
    @{ Cell = 'B30'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D30'; Value = '4.821' }
    @{ Cell = 'E30'; Value = '  -0.53%  ' }
    @{ Cell = 'B31'; Value = 'Stellar' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Cell = 'D31'; Value = '0.08715' }
    @{ Cell = 'E31'; Value = '  +0.23%  ' }
    @{ Cell = 'B32'; Value = 'HuobiToken' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D32'; Value = '3.102' }
    @{ Cell = 'E32'; Value = '  +2.28%  ' }
    @{ Cell = 'B33'; Value = 'ImmutableX' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D33'; Value = '0.7413' }
    @{ Cell = 'E33'; Value = '  +1.08%  ' }
    @{ Cell = 'B34'; Value = 'Filecoin' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Cell = 'D34'; Value = '4.425' }
    @{ Cell = 'E34'; Value = '  -1.90%  ' }
    @{ Cell = 'B35'; Value = 'RenderToken' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D35'; Value = '2.714' }
    @{ Cell = 'E35'; Value = '  +1.46%  ' }
    @{ Cell = 'B36'; Value = 'ARBITRUM' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D36'; Value = '1.110' }
    @{ Cell = 'E36'; Value = '  -0.98%  ' }
    @{ Cell = 'B37'; Value = 'TrustWalletToken' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D37'; Value = '1.069' }
    @{ Cell = 'E37'; Value = '  -1.41%  ' }
    @{ Cell = 'B38'; Value = 'VeChain' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D38'; Value = '0.01934' }
    @{ Cell = 'E38'; Value = '  -1.21%  ' }
    @{ Cell = 'B39'; Value = 'MXToken' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D39'; Value = '2.916' }
    @{ Cell = 'E39'; Value = '  +0.68%  ' }
    @{ Cell = 'B40'; Value = 'Hedera' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D40'; Value = '0.05078' }
    @{ Cell = 'E40'; Value = '  -0.78%  ' }
    @{ Cell = 'B41'; Value = 'TheSandbox' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Cell = 'D41'; Value = '0.5078' }
    @{ Cell = 'E41'; Value = '  +1.67%  ' }
    @{ Cell = 'B42'; Value = 'FraxShare' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D42'; Value = '6.779' }
    @{ Cell = 'E42'; Value = '  -3.04%  ' }
    @{ Cell = 'B43'; Value = 'Algorand' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = 'D43'; Value = '0.1507' }
    @{ Cell = 'E43'; Value = '  -3.09%  ' }
    @{ Cell = 'B44'; Value = 'Aptos' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D44'; Value = '8.020' }
    @{ Cell = 'E44'; Value = '  -1.81%  ' }
    @{ Cell = 'B45'; Value = 'Decentraland' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D45'; Value = '0.4680' }
    @{ Cell = 'E45'; Value = '  +1.58%  ' }
    @{ Cell = 'B46'; Value = 'PaxDollar' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Cell = 'D46'; Value = '1.003' }
    @{ Cell = 'E46'; Value = '  -0.56%  ' }
    @{ Cell = 'B47'; Value = 'EnergySwap' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D47'; Value = '9.992' }
    @{ Cell = 'E47'; Value = '  -0.35%  ' }
    @{ Cell = 'B48'; Value = 'Quant' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = 'D48'; Value = '99.06' }
    @{ Cell = 'E48'; Value = '  -2.20%  ' }
    @{ Cell = 'B49'; Value = 'NEARProtocol' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = 'D49'; Value = '1.566' }
    @{ Cell = 'E49'; Value = '  -1.51%  ' }
    @{ Cell = 'B50'; Value = 'Cronos' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D50'; Value = '0.05994' }
    @{ Cell = 'E50'; Value = '  -0.09%  ' }
    @{ Cell = 'B51'; Value = 'Aave' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D51'; Value = '63.60' }
    @{ Cell = 'E51'; Value = '  -1.30%  ' }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    # Leading apostrophe forces text storage so numeric-looking price
    # strings (e.g. '305.67', '19.30', '26.527.85') keep their exact
    # text instead of becoming floating point numbers.
    $cell.Value = "'" + $update.Value
    $cell.Style = "Normal"
}
